# Corrige (pone a 0) los valores de "Diferencia Stock" (columna L) para las
# filas indicadas, y actualiza el total "Total_Ajuste_Stock:" (C235) acorde.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(6,21,23,34,35,36,38,47,48,56,58,65,74,83,87,106,108,124,126,127,128,129,130,131,134,138,140,142,146,148,149,152,154,158,159,160,163,166,171,174,175,177,190,195,196,197,203,211,212,214,215)

foreach ($r in $rows) {
    $ws.Range("L$r").Value = 0
}

$ws.Range("C235").Value = 0
